# Update "想去人数" (number of people interested) counts in column F
# for the "展览" and "全部类型" worksheets.
#
# Row -> New value
#   F2  -> 257
#   F4  -> 288
#   F6  -> 289
#   F7  -> 6974
#   F8  -> 64
#   F11 -> 92
#   F12 -> 1
#   F18 -> 629
#   F19 -> 11

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 257
    "F4"  = 288
    "F6"  = 289
    "F7"  = 6974
    "F8"  = 64
    "F11" = 92
    "F12" = 1
    "F18" = 629
    "F19" = 11
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
